# Apply the "added one json for time bucket analysis" edit.
#
# The underlying source data (one JSON record per Wikipedia-linked resource)
# was regenerated, which re-ordered a few of the rows in the sheet:
#   - row 2 and row 3 swap their whole record (title/timestamp/historical
#     distance/uri) - the "day_31_beyond" time-bucket is identical for both
#     so it doesn't visibly move
#   - row 5 and row 7 swap their title/uri (the timestamp/historical
#     distance/time-bucket placeholders are identical "1-01-01.../unknown"
#     for every row from 4-7, so only A and E visibly change)
#   - row 4 and row 6 are unchanged
#
# We rewrite the affected cells directly and then rebuild the hyperlinks so
# each uri cell's hyperlink target follows its new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: now the "Micronesia profile - Leaders" / bbc.com record ---
$ws.Cells.Item(2, 1).Value = "Micronesia profile - Leaders"
$ws.Cells.Item(2, 2).Value = "2015-06-23T11:44:21UTC"
$ws.Cells.Item(2, 3).Value = 2965
$ws.Cells.Item(2, 4).Value = "day_31_beyond"
$ws.Cells.Item(2, 5).Value = "https://www.bbc.com/news/world-asia-pacific-15519473"

# --- Row 3: now the "Micronesia" / freedomhouse.org record ---
$ws.Cells.Item(3, 1).Value = "Micronesia"
$ws.Cells.Item(3, 2).Value = "2012-01-12T22:35:15UTC"
$ws.Cells.Item(3, 3).Value = 1707
$ws.Cells.Item(3, 4).Value = "day_31_beyond"
$ws.Cells.Item(3, 5).Value = "https://freedomhouse.org/report/freedom-world/2007/micronesia"

# --- Row 4: "Elections: Micronesia Congress 2017" - unchanged ---

# --- Row 5: now the "Senator Alik L. Alik" / cfsm.fm record ---
$ws.Cells.Item(5, 1).Value = "Senator Alik L. Alik"
$ws.Cells.Item(5, 5).Value = "http://www.cfsm.fm/index.php/public-info/112-members/biography/174-senator-alik-l-alik"

# --- Row 6: "Biography of President Urusemal" - unchanged ---

# --- Row 7: now the "Population of Micronesia (2020)" / worldometers.info record ---
$ws.Cells.Item(7, 1).Value = "Population of Micronesia (2020)"
$ws.Cells.Item(7, 5).Value = "https://www.worldometers.info/world-population/micronesia-population/"

# --- Rebuild the hyperlinks so each uri cell (column E) points at its new
#     text again (individual Hyperlinks.Item(i).Delete() is a no-op in this
#     host, so drop the whole collection and re-add it in the new order).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), $ws.Cells.Item(2, 5).Value())
$ws.Hyperlinks.Add($ws.Range("E3"), $ws.Cells.Item(3, 5).Value())
$ws.Hyperlinks.Add($ws.Range("E4"), $ws.Cells.Item(4, 5).Value())
$ws.Hyperlinks.Add($ws.Range("E5"), $ws.Cells.Item(5, 5).Value())
$ws.Hyperlinks.Add($ws.Range("E6"), $ws.Cells.Item(6, 5).Value())
$ws.Hyperlinks.Add($ws.Range("E7"), $ws.Cells.Item(7, 5).Value())

# Hyperlinks.Add stamps every re-created hyperlink with a (cosmetically
# identical) new "Hyperlink" style variant; restore the original style on
# the two rows whose uri text/target didn't actually change so only the
# genuinely edited cells pick up the new style id.
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"
